# Apply odds updates to Sheet1 as described by the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("G3").Value = 1.38
$ws.Range("H3").Value = 4.15
$ws.Range("I3").Value = 8.75
$ws.Range("J3").Value = 1.87
$ws.Range("K3").Value = 2.25
$ws.Range("L3").Value = 7.7
$ws.Range("M3").Value = 1.26
$ws.Range("N3").Value = 3.5
$ws.Range("O3").Value = 1.75
$ws.Range("P3").Value = 1.95
$ws.Range("Q3").Value = 2.77
$ws.Range("R3").Value = 1.39
$ws.Range("S3").Value = 1.39
$ws.Range("T3").Value = 2.77
$ws.Range("U3").Value = 2
$ws.Range("V3").Value = 1.72
$ws.Range("W3").Value = 6.4
$ws.Range("X3").Value = 6.3
$ws.Range("Z3").Value = 9
$ws.Range("AA3").Value = 11.25
$ws.Range("AC3").Value = 7.6
$ws.Range("AD3").Value = 8.25
$ws.Range("AE3").Value = 20
$ws.Range("AH3").Value = 20
$ws.Range("AI3").Value = 60
$ws.Range("AJ3").Value = 26
$ws.Range("AK3").Value = 250
$ws.Range("AL3").Value = 120
$ws.Range("AN3").Value = 1.05
$ws.Range("AO3").Value = 7.6

# Row 4
$ws.Range("G4").Value = 2.27
$ws.Range("H4").Value = 2.95
$ws.Range("I4").Value = 3.3
$ws.Range("J4").Value = 2.92
$ws.Range("K4").Value = 1.98
$ws.Range("M4").Value = 1.4
$ws.Range("N4").Value = 2.7
$ws.Range("O4").Value = 2.2
$ws.Range("P4").Value = 1.6
$ws.Range("Q4").Value = 3.8
$ws.Range("R4").Value = 1.23
$ws.Range("U4").Value = 1.87
$ws.Range("V4").Value = 1.83
$ws.Range("W4").Value = 6.6
$ws.Range("X4").Value = 10.25
$ws.Range("Z4").Value = 23
$ws.Range("AA4").Value = 20
$ws.Range("AB4").Value = 32
$ws.Range("AC4").Value = 6.1
$ws.Range("AD4").Value = 5.7
$ws.Range("AH4").Value = 8.5
$ws.Range("AI4").Value = 16.5
$ws.Range("AK4").Value = 50
$ws.Range("AM4").Value = 40
$ws.Range("AN4").Value = 1.09
$ws.Range("AO4").Value = 6.1

# Row 5
$ws.Range("G5").Value = 1.65
$ws.Range("H5").Value = 3.6
$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 2.25
$ws.Range("K5").Value = 2.12
$ws.Range("L5").Value = 5.2
$ws.Range("M5").Value = 1.28
$ws.Range("N5").Value = 3.35
$ws.Range("O5").Value = 1.82
$ws.Range("P5").Value = 1.88
$ws.Range("Q5").Value = 2.95
$ws.Range("R5").Value = 1.35
$ws.Range("T5").Value = 2.67
$ws.Range("U5").Value = 1.8
$ws.Range("V5").Value = 1.91
$ws.Range("X5").Value = 7.9
$ws.Range("Y5").Value = 7.9
$ws.Range("Z5").Value = 13
$ws.Range("AA5").Value = 13
$ws.Range("AC5").Value = 7.5
$ws.Range("AD5").Value = 7
$ws.Range("AH5").Value = 13.5
$ws.Range("AI5").Value = 30
$ws.Range("AJ5").Value = 15.5
$ws.Range("AK5").Value = 90
$ws.Range("AL5").Value = 50
$ws.Range("AM5").Value = 50
$ws.Range("AO5").Value = 7.5

# Row 13
$ws.Range("L13").Value = 4.2
$ws.Range("Q13").Value = 2.7
$ws.Range("R13").Value = 1.35
$ws.Range("AB13").Value = 25
$ws.Range("AM13").Value = 40

# Row 14
$ws.Range("G14").Value = 1.57
$ws.Range("U14").Value = 2.38
$ws.Range("V14").Value = 1.53
$ws.Range("X14").Value = 6.5
$ws.Range("Z14").Value = 11
$ws.Range("AI14").Value = 26

# Row 17
$ws.Range("G17").Value = 3
$ws.Range("I17").Value = 2.9
$ws.Range("J17").Value = 4
$ws.Range("L17").Value = 4
$ws.Range("W17").Value = 6
$ws.Range("X17").Value = 12
$ws.Range("Y17").Value = 13
$ws.Range("Z17").Value = 34
$ws.Range("AA17").Value = 34
$ws.Range("AD17").Value = 6
$ws.Range("AI17").Value = 12
$ws.Range("AJ17").Value = 13
$ws.Range("AK17").Value = 34
